$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A89").Value = "2025-04-29 14:53:09"
$ws.Range("B89").Value = 244
